# Update the test/training schedule data in row 3 (trial 2):
#   E3 (y_corrSteps): 5 -> 4
#   G3 (y_nrSteps):  -2 -> -3
#   H3 (alienID):    14 -> 13
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 4
$ws.Range("G3").Value = -3
$ws.Range("H3").Value = 13

# Leave the selection on E3, matching the cell that was edited.
$ws.Range("E3").Select() | Out-Null
